$d = $word.ActiveDocument

# Locate the bullet paragraph to remove: "Clean up memory reading/writing API ..."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Clean up memory reading/writing API*") {
        $target = $p
        break
    }
}

$delStart = $target.Range.Start

# Delete the whole paragraph (text + paragraph mark). The following
# paragraph's runs merge into this <w:p>, keeping its pPr - exactly what
# Word does when you select-and-delete an entire bullet line.
$target.Range.Delete()

# Word leaves its "last edit" marker (the hidden _GoBack bookmark) at the
# point of the edit, so relocate it from wherever it used to be to here.
$here = $d.Range($delStart, $delStart)
$d.Bookmarks.Add("_GoBack", $here)

# The old _GoBack location (end of the final "Support 64-bit parameters ..."
# bullet) used to split that bullet's trailing text into two runs, one
# holding just the final "."; with the bookmark gone, re-touch that "."
# in place so it's merged back into the run before it.
$lastBullet = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Support 64-bit parameters under x86 in remote function caller.*") {
        $lastBullet = $p
    }
}
$endRange = $lastBullet.Range
$dotPos = $endRange.End - 2
$dotRange = $d.Range($dotPos, $dotPos + 1)
$dotRange.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)
